$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: same bold/border/centered style as the other
# header cells (copy formatting from G1), then set its text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New data cell H2 (plain, unstyled numeric 0 — matches the "Save" column).
$ws.Range("H2").Value = 0
